$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting the existing rows 10 and 11
# down to rows 11 and 12 respectively.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(10, 3).Value = "Los Lagos"
$ws.Cells.Item(10, 4).Value = 44474
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(10, 6).Value = 300000000
$ws.Cells.Item(10, 7).Value = "Espárragos"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 250
$ws.Cells.Item(10, 11).Value = 2000
$ws.Cells.Item(10, 12).Value = 2000
$ws.Cells.Item(10, 13).Value = 2000
$ws.Cells.Item(10, 14).Value = "`$/kilo"
$ws.Cells.Item(10, 15).Value = "Provincia de Linares"
$ws.Cells.Item(10, 16).Value = 2000
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# Match the date formatting used by the other "Fecha" cells in column D.
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat
